$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the last-updated timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 23:08"

# Azerbaiyan overtakes Kenia in ranking (rows 73/74 swap countries + refreshed counts)
$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("A74").Value = "Kenia"

# Trinidad y Tobago overtakes Congo in ranking (rows 130/131 swap countries + refreshed counts)
$ws.Range("A130").Value = "Trinidad yTobago"
$ws.Range("A131").Value = "Congo"

# Refresh updated case-count figures
$ws.Range("B4").Value = 8206034
$ws.Range("C4").Value = 55991
$ws.Range("D4").Value = 5308029
$ws.Range("E4").Value = 2675372
$ws.Range("G4").Value = 790
$ws.Range("H4").Value = 222633
$ws.Range("B22").Value = 348816
$ws.Range("C22").Value = 7074
$ws.Range("E22").Value = 57106
$ws.Range("B27").Value = 300201
$ws.Range("C27").Value = 1701
$ws.Range("D27").Value = 257226
$ws.Range("E27").Value = 40848
$ws.Range("G27").Value = 29
$ws.Range("H27").Value = 2127
$ws.Range("B30").Value = 191344
$ws.Range("C30").Value = 1957
$ws.Range("D30").Value = 161151
$ws.Range("E30").Value = 20495
$ws.Range("B50").Value = 93152
$ws.Range("C50").Value = 1372
$ws.Range("D50").Value = 57731
$ws.Range("E50").Value = 34262
$ws.Range("G50").Value = 25
$ws.Range("H50").Value = 1159
$ws.Range("B73").Value = 43280
$ws.Range("C73").Value = 530
$ws.Range("D73").Value = 39671
$ws.Range("E73").Value = 2990
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 619
$ws.Range("B74").Value = 43143
$ws.Range("C74").Value = 602
$ws.Range("D74").Value = 31508
$ws.Range("E74").Value = 10830
$ws.Range("G74").Value = 8
$ws.Range("H74").Value = 805
$ws.Range("B92").Value = 20257
$ws.Range("C92").Value = 40
$ws.Range("D92").Value = 19898
$ws.Range("E92").Value = 239
$ws.Range("B103").Value = 12103
$ws.Range("C103").Value = 34
$ws.Range("D103").Value = 10329
$ws.Range("E103").Value = 1644
$ws.Range("B105").Value = 11113
$ws.Range("C105").Value = 51
$ws.Range("D105").Value = 9931
$ws.Range("B124").Value = 5733
$ws.Range("C124").Value = 18
$ws.Range("D124").Value = 5375
$ws.Range("E124").Value = 243
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 115
$ws.Range("B125").Value = 5443
$ws.Range("C125").Value = 3
$ws.Range("D125").Value = 5368
$ws.Range("E125").Value = 14
$ws.Range("B130").Value = 5194
$ws.Range("C130").Value = 40
$ws.Range("D130").Value = 3502
$ws.Range("E130").Value = 1599
$ws.Range("H130").Value = 93
$ws.Range("B131").Value = 5156
$ws.Range("D131").Value = 3887
$ws.Range("E131").Value = 1177
$ws.Range("H131").Value = 92
$ws.Range("B134").Value = 4953
$ws.Range("C134").Value = 13
$ws.Range("D134").Value = 4601
$ws.Range("E134").Value = 319
$ws.Range("B146").Value = 3644
$ws.Range("C146").Value = 2
$ws.Range("D146").Value = 2646
$ws.Range("E146").Value = 880
$ws.Range("B157").Value = 2323
$ws.Range("C157").Value = 8
$ws.Range("D157").Value = 1746
$ws.Range("E157").Value = 504
$ws.Range("B159").Value = 2285
$ws.Range("C159").Value = 104
$ws.Range("E159").Value = 816
$ws.Range("B165").Value = 1374
$ws.Range("C165").Value = 2
$ws.Range("E165").Value = 38
$ws.Range("B166").Value = 1350
$ws.Range("C166").Value = 21
$ws.Range("E166").Value = 138
$ws.Range("B174").Value = 673
$ws.Range("C174").Value = 28
$ws.Range("D174").Value = 371
$ws.Range("E174").Value = 301
